$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 862 (shifts existing rows 862:903 down to 863:904,
# and Excel auto-extends the used range / dimension accordingly).
$ws.Rows.Item(862).Insert()

# Fill the newly inserted row with the new data point for 2026/02/26 (Thu).
# Column A/B hold text (date-like / weekday strings) in the source data, so
# force Text formatting before assigning to avoid Excel auto-converting the
# "yyyy/mm/dd" string into a date serial number, then clear the formatting
# again so the cell ends up with no explicit style (matching the rest of
# the data rows).
$ws.Range("A862").NumberFormat = "@"
$ws.Range("A862").Value = "2026/02/26"
$ws.Range("A862").ClearFormats()

$ws.Range("B862").NumberFormat = "@"
$ws.Range("B862").Value = "木"
$ws.Range("B862").ClearFormats()

$ws.Range("C862").Value = 14
$ws.Range("D862").Value = 21
